$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.287.55"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.865.02"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'236.44"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4716"
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("D8").Value = "'0.2903"
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").Value = "'0.06545"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "'21.90"
$ws.Range("E10").Value = "  +3.29%  "
$ws.Range("D11").Value = "'0.07944"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "'98.02"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "1.863.04"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "'5.152"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "'0.6810"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "'263.60"
$ws.Range("E16").Value = "  -6.28%  "
$ws.Range("D17").Value = "30.272.15"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'13.77"
$ws.Range("E18").Value = "  +8.40%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'0.000007466"
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").Value = "2.107.24"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'5.267"
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("D24").Value = "'6.171"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'167.43"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "'9.192"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "'18.91"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "'1.951"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "'1.396"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("D30").Value = "'0.09859"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").Value = "'4.351"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").Value = "'1.473"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "'4.025"
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("D34").Value = "'0.04718"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").Value = "'1.131"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "'0.7005"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "'2.709"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "'0.01879"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'2.628"
$ws.Range("E39").Value = "  +3.42%  "
$ws.Range("D40").Value = "'6.338"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "'73.89"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "'1.943"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").Value = "'0.8436"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'0.4158"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "'103.33"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'7.159"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").Value = "'949.41"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "'9.206"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'34.15"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "'0.05664"
$ws.Range("E51").Value = "  +0.54%  "
